$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row=2; D="66.869.34"; E="  +3.37%  " }
    @{ Row=3; D="3.094.82"; E="  +5.51%  " }
    @{ Row=4; D=$null; E="  -0.05%  " }
    @{ Row=5; D="580.24"; E="  +2.42%  " }
    @{ Row=6; D="167.95"; E="  +6.93%  " }
    @{ Row=7; D=$null; E="  -0.10%  " }
    @{ Row=8; D="3.089.35"; E="  +5.50%  " }
    @{ Row=9; D=$null; E="  +1.66%  " }
    @{ Row=10; D="6.67"; E="  +0.58%  " }
    @{ Row=11; D=$null; E="  +3.19%  " }
    @{ Row=12; D="0.483"; E="  +6.16%  " }
    @{ Row=13; D="0.0000250"; E="  +2.70%  " }
    @{ Row=14; D="36.72"; E="  +8.54%  " }
    @{ Row=15; D=$null; E="  -0.72%  " }
    @{ Row=16; D="3.608.64"; E="  +5.31%  " }
    @{ Row=17; D="66.880.43"; E="  +2.88%  " }
    @{ Row=18; D="7.24"; E="  +4.77%  " }
    @{ Row=19; D="3.096.60"; E="  +5.39%  " }
    @{ Row=20; D="16.32"; E="  +18.34%  " }
    @{ Row=21; D="469.45"; E="  +5.47%  " }
    @{ Row=22; D=$null; E="  +5.61%  " }
    @{ Row=23; D="7.57"; E="  +5.03%  " }
    @{ Row=24; D="83.29"; E="  +0.92%  " }
    @{ Row=25; D=$null; E="  +8.80%  " }
    @{ Row=26; D="12.87"; E="  +7.23%  " }
    @{ Row=27; D=$null; E="  +1.87%  " }
    @{ Row=28; D=$null; E="  +0.07%  " }
    @{ Row=29; D=$null; E="  +2.28%  " }
    @{ Row=30; D=$null; E="  +3.52%  " }
    @{ Row=31; D=$null; E="  +4.37%  " }
    @{ Row=32; D=$null; E="  +3.93%  " }
    @{ Row=33; D=$null; E="  +4.47%  " }
    @{ Row=34; D=$null; E="  +5.56%  " }
    @{ Row=35; D=$null; E="  +0.07%  " }
    @{ Row=36; D="1.00"; E="  +3.16%  " }
    @{ Row=37; D=$null; E="  +4.34%  " }
    @{ Row=38; D="46.74"; E="  +7.88%  " }
    @{ Row=39; D=$null; E="  +5.99%  " }
    @{ Row=40; D="50.32"; E="  +2.60%  " }
    @{ Row=41; D="0.317"; E="  +6.99%  " }
    @{ Row=42; D="0.123"; E="  +4.40%  " }
    @{ Row=43; D="8.71"; E="  +3.80%  " }
    @{ Row=44; D=$null; E="  +2.60%  " }
    @{ Row=45; D="0.0363"; E="  +3.69%  " }
    @{ Row=46; D="387.53"; E="  +1.88%  " }
    @{ Row=47; D="2.772.57"; E="  +1.93%  " }
    @{ Row=48; D="134.78"; E="  +2.29%  " }
    @{ Row=49; D=$null; E="  -0.01%  " }
    @{ Row=50; D="24.84"; E="  +7.55%  " }
    @{ Row=51; D=$null; E="  +5.42%  " }
)

foreach ($item in $changes) {
    $row = $item.Row

    if ($null -ne $item.D) {
        # Prefix with an apostrophe so values that look numeric (e.g. "1.00",
        # "0.0000250") are stored as literal text, preserving exact digits
        # instead of being auto-coerced into a Number by Excel.
        $ws.Cells.Item($row, 4).Value = "'" + $item.D
    }

    if ($null -ne $item.E) {
        $ws.Cells.Item($row, 5).Value = $item.E
    }
}
